# Apply updated crypto price/volume snapshot values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text storage for numeric-looking prices (matches
# the source data's inline-string / text cell type instead of Excel auto-
# converting them to numbers, which would lose trailing zeros / formatting).
$textPrefix = [string][char]39

$ws.Range('D2').Value = '55.182.71'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '2.346.29'
$ws.Range('E3').Value = '  -4.34%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = $textPrefix + '476.12'
$ws.Range('E5').Value = '  -1.78%  '
$ws.Range('D6').Value = $textPrefix + '145.76'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('D7').Value = $textPrefix + '0.999'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = $textPrefix + '0.610'
$ws.Range('E8').Value = '  +20.75%  '
$ws.Range('D9').Value = '2.357.59'
$ws.Range('E9').Value = '  -4.18%  '
$ws.Range('D10').Value = $textPrefix + '0.0959'
$ws.Range('E10').Value = '  -0.25%  '
$ws.Range('E11').Value = '  -5.49%  '
$ws.Range('D12').Value = $textPrefix + '0.325'
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').Value = '2.758.41'
$ws.Range('E14').Value = '  -4.32%  '
$ws.Range('D15').Value = '55.163.63'
$ws.Range('E15').Value = '  -1.38%  '
$ws.Range('D16').Value = $textPrefix + '19.94'
$ws.Range('E16').Value = '  -4.42%  '
$ws.Range('D17').Value = $textPrefix + '0.0000129'
$ws.Range('E17').Value = '  -4.08%  '
$ws.Range('D18').Value = '2.357.91'
$ws.Range('E18').Value = '  -4.31%  '
$ws.Range('D19').Value = $textPrefix + '4.58'
$ws.Range('E19').Value = '  +2.06%  '
$ws.Range('D20').Value = $textPrefix + '314.14'
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('D21').Value = $textPrefix + '9.55'
$ws.Range('E21').Value = '  -4.55%  '
$ws.Range('D22').Value = $textPrefix + '0.998'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').Value = $textPrefix + '5.67'
$ws.Range('D24').Value = $textPrefix + '56.38'
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('D26').Value = $textPrefix + '0.394'
$ws.Range('E26').Value = '  -3.42%  '
$ws.Range('E27').Value = '  -3.93%  '
$ws.Range('D28').Value = '2.453.45'
$ws.Range('E28').Value = '  -4.66%  '
$ws.Range('D29').Value = $textPrefix + '7.05'
$ws.Range('E29').Value = '  -8.09%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').Value = '0.0₃0742'
$ws.Range('E31').Value = '  -4.30%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = $textPrefix + '145.74'
$ws.Range('E32').Value = '  -1.19%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = $textPrefix + '18.12'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('D35').Value = $textPrefix + '5.10'
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = $textPrefix + '1.09'
$ws.Range('E36').Value = '  -3.93%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = $textPrefix + '3.59'
$ws.Range('E37').Value = '  -2.67%  '
$ws.Range('D38').Value = $textPrefix + '0.806'
$ws.Range('E38').Value = '  -5.06%  '
$ws.Range('E39').Value = '  +10.80%  '
$ws.Range('D40').Value = $textPrefix + '33.68'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = $textPrefix + '0.998'
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = $textPrefix + '3.37'
$ws.Range('E42').Value = '  -3.61%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = $textPrefix + '1.31'
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('E44').Value = '  -3.44%  '
$ws.Range('D45').Value = $textPrefix + '0.0516'
$ws.Range('E45').Value = '  -5.66%  '
$ws.Range('D46').Value = $textPrefix + '10.16'
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = $textPrefix + '0.0220'
$ws.Range('E47').Value = '  -2.45%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = $textPrefix + '248.50'
$ws.Range('E48').Value = '  -4.36%  '
$ws.Range('D49').Value = $textPrefix + '4.38'
$ws.Range('E49').Value = '  -6.25%  '
$ws.Range('D50').Value = '1.798.12'
$ws.Range('D51').Value = $textPrefix + '16.61'
$ws.Range('E51').Value = '  -4.55%  '
